# Applies the cryptos-list refresh described by the commit
# "Updated cryptos list ... with GitHub Actions": new Price (col D)
# and Volume(1h) (col E) figures for every coin row, plus a swap of
# the Fetch.AI / PancakeSwap rows (29/30 by 0-index -> rows 31/32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per changed record. Only the columns that actually change
# are present in each hashtable (B/C for the two swapped coins,
# D/E -- or just E -- for the rest).
$updates = @(
    @{Row=2; D="57.982.62"; E="  -0.37%  "},
    @{Row=3; D="2.466.39"; E="  -0.50%  "},
    @{Row=4; D="0.999"; E="  -0.21%  "},
    @{Row=5; D="513.64"; E="  -1.27%  "},
    @{Row=6; D="130.31"; E="  -0.75%  "},
    @{Row=7; D="0.996"; E="  -0.42%  "},
    @{Row=8; D="0.551"; E="  -1.48%  "},
    @{Row=9; D="2.488.77"; E="  +0.21%  "},
    @{Row=10; D="0.0965"; E="  -2.97%  "},
    @{Row=11; D="0.157"; E="  -0.07%  "},
    @{Row=12; D="5.20"; E="  -2.96%  "},
    @{Row=13; D="0.329"; E="  -4.33%  "},
    @{Row=14; D="2.899.94"; E="  -0.74%  "},
    @{Row=15; D="57.886.38"; E="  -0.42%  "},
    @{Row=16; D="21.95"; E="  -1.94%  "},
    @{Row=17; D="0.0000134"; E="  -1.93%  "},
    @{Row=18; D="2.472.38"; E="  -0.35%  "},
    @{Row=19; D="10.57"; E="  -2.74%  "},
    @{Row=20; D="318.75"; E="  -0.63%  "},
    @{Row=21; D="4.13"; E="  -1.41%  "},
    @{Row=22; D="1.00"; E="  +0.11%  "},
    @{Row=23; D="5.94"; E="  +2.91%  "},
    @{Row=24; D="62.99"; E="  -1.88%  "},
    @{Row=25; D="0.401"; E="  -2.35%  "},
    @{Row=26; D="0.991"; E="  -0.91%  "},
    @{Row=27; E="  -0.09%  "},
    @{Row=28; D="7.26"; E="  -1.32%  "},
    @{Row=29; D="169.43"; E="  +2.03%  "},
    @{Row=30; D="0.0₃0733"; E="  -3.15%  "},
    @{Row=31; B="Fetch.AI"; C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D="1.17"; E="  +0.03%  "},
    @{Row=32; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="1.67"; E="  -2.01%  "},
    @{Row=33; D="6.21"; E="  -1.96%  "},
    @{Row=34; D="0.997"; E="  -0.15%  "},
    @{Row=35; E="  -0.52%  "},
    @{Row=36; D="17.82"; E="  -1.67%  "},
    @{Row=37; D="1.27"; E="  -3.88%  "},
    @{Row=38; D="3.91"; E="  -2.02%  "},
    @{Row=39; D="36.63"; E="  +0.08%  "},
    @{Row=40; D="1.46"; E="  -1.57%  "},
    @{Row=41; D="0.766"; E="  -3.18%  "},
    @{Row=42; D="272.27"; E="  -1.61%  "},
    @{Row=43; D="5.04"; E="  +0.22%  "},
    @{Row=44; D="3.38"; E="  -3.04%  "},
    @{Row=45; D="0.588"; E="  -1.14%  "},
    @{Row=46; E="  +0.73%  "},
    @{Row=47; D="121.35"; E="  -4.83%  "},
    @{Row=48; D="0.0488"; E="  -0.28%  "},
    @{Row=49; D="17.45"; E="  -3.01%  "},
    @{Row=50; D="0.0211"; E="  -1.64%  "},
    @{Row=51; D="16.69"; E="  -2.60%  "}
)

# Column D holds plain-decimal-looking price strings (e.g. "0.999",
# "1.00", "513.64") that Excel would otherwise auto-convert to the
# Number type on assignment, silently dropping significant trailing
# zeros ("1.00" -> 1). Force column D to Text for the duration of the
# write, then restore the default "Normal" style so the sheet keeps
# its original (unstyled) look, matching every other data cell.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Cells.Item($u.Row, 4).Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}

$priceRange.Style = "Normal"

Write-Host "Applied $($updates.Count) row updates."
